# Apply the changes described by the diff:
#  1. Add a new row value "נחום" in cell D7 (new shared string).
#  2. Move the active selection from D5 to D6.
#  3. Nudge the saved workbook window position (xWindow) from 930 to 2050.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# 1. Set the new cell value - this also appends "נחום" to the shared strings table.
$ws.Range("D7").Value = "נחום"

# 2. Update the selected/active cell shown in the worksheet view.
$ws.Range("D6").Select()

# 3. Record the updated window position for the workbook view.
$win = $excel.ActiveWindow
$win.Left = 2050
